$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Chad"
$ws.Cells.Item(2,3).Value = "Itga2"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.6240709999999999
$ws.Cells.Item(2,8).Value = 1.872213
$ws.Cells.Item(2,9).Value = 0.07908544873752882
$ws.Cells.Item(2,10).Value = 0.07908544873752882
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 0.229822
$ws.Cells.Item(2,14).Value = 0.689466
$ws.Cells.Item(2,15).Value = 0.09226175421862418
$ws.Cells.Item(2,16).Value = 0.09226175421862419
$ws.Cells.Item(2,17).Value = 0.143425245362
$ws.Cells.Item(2,18).Value = 1.290827208258
$ws.Cells.Item(2,19).Value = 0.007296562233691486
$ws.Cells.Item(2,20).Value = 0.007296562233691487

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Chad"
$ws.Cells.Item(3,3).Value = "Itga2"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.6240709999999999
$ws.Cells.Item(3,8).Value = 1.872213
$ws.Cells.Item(3,9).Value = 0.07908544873752882
$ws.Cells.Item(3,10).Value = 0.07908544873752882
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 1.087098333333333
$ws.Cells.Item(3,14).Value = 3.261295
$ws.Cells.Item(3,15).Value = 0.4364142651333466
$ws.Cells.Item(3,16).Value = 0.4364142651333466
$ws.Cells.Item(3,17).Value = 0.6784265439816666
$ws.Cells.Item(3,18).Value = 6.105838895834999
$ws.Cells.Item(3,19).Value = 0.03451401799352959
$ws.Cells.Item(3,20).Value = 0.03451401799352959

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Chad"
$ws.Cells.Item(4,3).Value = "Itga2"
$ws.Cells.Item(4,4).Value = "MuSCs"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.6240709999999999
$ws.Cells.Item(4,8).Value = 1.872213
$ws.Cells.Item(4,9).Value = 0.07908544873752882
$ws.Cells.Item(4,10).Value = 0.07908544873752882
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 1.174057666666666
$ws.Cells.Item(4,14).Value = 3.522173
$ws.Cells.Item(4,15).Value = 0.4713239806480292
$ws.Cells.Item(4,16).Value = 0.4713239806480293
$ws.Cells.Item(4,17).Value = 0.7326953420943331
$ws.Cells.Item(4,18).Value = 6.594258078848998
$ws.Cells.Item(4,19).Value = 0.03727486851030774
$ws.Cells.Item(4,20).Value = 0.03727486851030775

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Chad"
$ws.Cells.Item(5,3).Value = "Itga2"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 6.439664333333334
$ws.Cells.Item(5,8).Value = 19.318993
$ws.Cells.Item(5,9).Value = 0.8160669916094901
$ws.Cells.Item(5,10).Value = 0.8160669916094901
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.229822
$ws.Cells.Item(5,14).Value = 0.689466
$ws.Cells.Item(5,15).Value = 0.09226175421862418
$ws.Cells.Item(5,16).Value = 0.09226175421862419
$ws.Cells.Item(5,17).Value = 1.479976536415333
$ws.Cells.Item(5,18).Value = 13.319788827738
$ws.Cells.Item(5,19).Value = 0.07529177220580682
$ws.Cells.Item(5,20).Value = 0.07529177220580682

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Chad"
$ws.Cells.Item(6,3).Value = "Itga2"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 6.439664333333334
$ws.Cells.Item(6,8).Value = 19.318993
$ws.Cells.Item(6,9).Value = 0.8160669916094901
$ws.Cells.Item(6,10).Value = 0.8160669916094901
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 1.087098333333333
$ws.Cells.Item(6,14).Value = 3.261295
$ws.Cells.Item(6,15).Value = 0.4364142651333466
$ws.Cells.Item(6,16).Value = 0.4364142651333466
$ws.Cells.Item(6,17).Value = 7.000548363992778
$ws.Cells.Item(6,18).Value = 63.004935275935
$ws.Cells.Item(6,19).Value = 0.3561432764428366
$ws.Cells.Item(6,20).Value = 0.3561432764428366

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Chad"
$ws.Cells.Item(7,3).Value = "Itga2"
$ws.Cells.Item(7,4).Value = "MuSCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 6.439664333333334
$ws.Cells.Item(7,8).Value = 19.318993
$ws.Cells.Item(7,9).Value = 0.8160669916094901
$ws.Cells.Item(7,10).Value = 0.8160669916094901
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 1.174057666666666
$ws.Cells.Item(7,14).Value = 3.522173
$ws.Cells.Item(7,15).Value = 0.4713239806480292
$ws.Cells.Item(7,16).Value = 0.4713239806480293
$ws.Cells.Item(7,17).Value = 7.560537281309888
$ws.Cells.Item(7,18).Value = 68.044835531789
$ws.Cells.Item(7,19).Value = 0.3846319429608467
$ws.Cells.Item(7,20).Value = 0.3846319429608469

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Chad"
$ws.Cells.Item(8,3).Value = "Itga2"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.8273623333333333
$ws.Cells.Item(8,8).Value = 2.482087
$ws.Cells.Item(8,9).Value = 0.1048475596529811
$ws.Cells.Item(8,10).Value = 0.1048475596529811
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.229822
$ws.Cells.Item(8,14).Value = 0.689466
$ws.Cells.Item(8,15).Value = 0.09226175421862418
$ws.Cells.Item(8,16).Value = 0.09226175421862419
$ws.Cells.Item(8,17).Value = 0.1901460661713333
$ws.Cells.Item(8,18).Value = 1.711314595542
$ws.Cells.Item(8,19).Value = 0.009673419779125877
$ws.Cells.Item(8,20).Value = 0.00967341977912588

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Chad"
$ws.Cells.Item(9,3).Value = "Itga2"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.8273623333333333
$ws.Cells.Item(9,8).Value = 2.482087
$ws.Cells.Item(9,9).Value = 0.1048475596529811
$ws.Cells.Item(9,10).Value = 0.1048475596529811
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 1.087098333333333
$ws.Cells.Item(9,14).Value = 3.261295
$ws.Cells.Item(9,15).Value = 0.4364142651333466
$ws.Cells.Item(9,16).Value = 0.4364142651333466
$ws.Cells.Item(9,17).Value = 0.8994242136294444
$ws.Cells.Item(9,18).Value = 8.094817922665
$ws.Cells.Item(9,19).Value = 0.04575697069698046
$ws.Cells.Item(9,20).Value = 0.04575697069698046

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Chad"
$ws.Cells.Item(10,3).Value = "Itga2"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.8273623333333333
$ws.Cells.Item(10,8).Value = 2.482087
$ws.Cells.Item(10,9).Value = 0.1048475596529811
$ws.Cells.Item(10,10).Value = 0.1048475596529811
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 1.174057666666666
$ws.Cells.Item(10,14).Value = 3.522173
$ws.Cells.Item(10,15).Value = 0.4713239806480292
$ws.Cells.Item(10,16).Value = 0.4713239806480293
$ws.Cells.Item(10,17).Value = 0.9713710905612221
$ws.Cells.Item(10,18).Value = 8.742339815050999
$ws.Cells.Item(10,19).Value = 0.04941716917687475
$ws.Cells.Item(10,20).Value = 0.04941716917687476
